$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2, pushing existing data rows down
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "schizophrenia_controls_sc_w2v"
$ws.Range("B2").Value = 0.8666666666666668
$ws.Range("C2").Value = 0.19436506316151
$ws.Range("D2").Value = "gaussian-nb"
